$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refactor "Custom EuTaxonomyAmountWithCurrencyComponent" -> "AmountWithCurrencyComponent"
# for all cells in column F that reference it. The other two "Custom EuTaxonomy*Component"
# shared strings keep their text; they only shift shared-string indices once the old
# string is no longer referenced.
$cellsToRename = @("F12", "F14", "F16", "F19", "F31", "F33", "F35", "F38", "F50", "F52", "F54", "F57")
foreach ($addr in $cellsToRename) {
    $ws.Range($addr).Value2 = "AmountWithCurrencyComponent"
}

# Update the sheet view's scroll position / active cell to match the latest edit session.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 46
$ws.Range("F57").Select()
